$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.250120759010315
$ws.Range("B1").Value = 2.640519857406616
$ws.Range("C1").Value = 8.353707313537598
$ws.Range("D1").Value = 2.106642246246338
$ws.Range("E1").Value = 1.13615083694458
